$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.105.15'
$ws.Range("E2").Value = '  +3.88%  '
$ws.Range("D3").Value = '2.324.20'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'545.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.70%  '
$ws.Range("D6").Value = "'130.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'0.577"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").Value = '2.320.95'
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("D11").Value = "'5.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = "'0.333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = "'23.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.64%  '
$ws.Range("D15").Value = '2.737.25'
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '60.083.19'
$ws.Range("E16").Value = '  +3.92%  '
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").Value = '2.325.41'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").Value = "'10.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").Value = "'4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").Value = "'313.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("D22").Value = "'6.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.14%  '
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("D24").Value = "'63.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("D25").Value = "'0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").Value = "'7.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("D28").Value = "'1.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.90%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = "'1.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.35%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = "'172.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +13.18%  '
$ws.Range("D32").Value = '0.0₃0724'
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("D33").Value = "'5.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.05%  '
$ws.Range("D34").Value = "'1.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.48%  '
$ws.Range("D35").Value = "'0.378"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").Value = "'17.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").Value = "'4.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.19%  '
$ws.Range("D40").Value = "'325.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.53%  '
$ws.Range("D41").Value = "'37.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("D44").Value = "'3.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("D45").Value = "'0.0943"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = "'19.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.27%  '
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("D48").Value = "'0.559"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0228'
$ws.Range("E49").Value = '  +23.17%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = "'0.0213"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.95%  '
$ws.Range("D51").Value = "'11.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.59%  '
